# Restore cell C10 on the "Rules" sheet from 18 to 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$ws.Cells.Item(10, 3).Value = 1
